# "Changes for new model" - rework both grade sheets to the new
# Matematica / Comunicación / Ingles layout.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename the sheets (tabs) to the new model's names.
$ws1.Name = "Primer Grado"
$ws2.Name = "Segundo Grado"

foreach ($ws in @($ws1, $ws2)) {
    # Drop the old Curso D/E/F columns (and any leftover values in them).
    $ws.Range("E1:G3").Clear() | Out-Null

    # New subject headers.
    $ws.Range("A1").Value = "Alumno"
    $ws.Range("B1").Value = "Matematica"
    $ws.Range("C1").Value = "Comunicación"
    $ws.Range("D1").Value = "Ingles"

    # New sample data rows.
    $ws.Range("A2").Value = "Mark"
    $ws.Range("B2").Value = 14
    $ws.Range("C2").Value = 15
    $ws.Range("D2").Value = 12

    $ws.Range("A3").Value = "Jose"
    $ws.Range("B3").Value = 15
    $ws.Range("C3").Value = 12
    $ws.Range("D3").Value = 12
}

# Primer Grado's header row now wraps its (longer) subject names.
$ws1.Range("A1:D1").WrapText = $true

# Give column C a bit more breathing room for "Comunicación" on both sheets.
$ws1.Columns.Item(3).ColumnWidth = 12.7
$ws2.Columns.Item(3).ColumnWidth = 13.5

Write-Host "Applied new-model changes"
